$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q16:Q23").ClearContents()
